$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = $null
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $true
$ws.Range("C4").Value = $null

$ws.Range("C1:C5").Select()
